# docs : WBS 업데이트
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2_WBS")

# Row 10 - 포스터: progress 0 -> 10
$ws.Range("I10").Value = 10

# Row 17 - 리소스 기획서: progress 55 -> 65
$ws.Range("I17").Value = 65

# Row 25 - "컨텐츠" task becomes "엔딩 일러스트", progress 0 -> 30
$ws.Range("C25").Value = "엔딩 일러스트"
$ws.Range("I25").Value = 30

# Row 26 - UI: progress 25 -> 45
$ws.Range("I26").Value = 45

# Rows 28-31 (3D 모델링 milestone block)
# Row 28 - 모델링: done (was tracked as 1% on a 0-1 percent-style format, now 100 on a literal "0%" custom format)
$ws.Range("I28").NumberFormat = "0""%"""
$ws.Range("I28").Value = 100

# Row 29 - 리깅: completed 2021-11-08, progress 100
$ws.Range("H29").Value = 44508
$ws.Range("I29").NumberFormat = "0""%"""
$ws.Range("I29").Value = 100

# Row 30 - 모션: completed 2021-11-15, progress 100
$ws.Range("H30").Value = 44515
$ws.Range("I30").NumberFormat = "0""%"""
$ws.Range("I30").Value = 100

# Row 31 - new task "텍스쳐" added to the milestone block: 2021-10-09 ~ 2021-11-17, completed 2021-11-08, progress 100
$ws.Range("C31").Value = "텍스쳐"
$ws.Range("F31").Value = 44478
$ws.Range("G31").Value = 44517
$ws.Range("H31").Value = 44508
$ws.Range("I31").NumberFormat = "0""%"""
$ws.Range("I31").Value = 100

$ws.Range("M13").Select()
